# -----------------------------------------------------------------------
# contrainte ajoutee : le joueur ne peut jouer qu'une seule fois a la
# question posee
#
# Two edits to the document body:
#
#  1. Paragraph "Les technologies utilisees... : html, css, javascript,
#     php" - the stray _GoBack bookmark sitting between "utilisees" and
#     " : html, " is removed and the two runs it used to separate are
#     merged back into one.
#
#  2. Paragraph "On a appris pas mal de chose...fonctionnalites..." -
#     the trailing clause describing the still-missing "can't answer
#     the same question twice" feature is replaced by "demandees"
#     (i.e. the feature is now considered done), and the _GoBack
#     bookmark is moved to sit right after the new word, at the very
#     end of the paragraph.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1 ----------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-match the exact (unchanged) text spanning the old run boundary and
# "replace" it with itself; this collapses the two adjacent runs into a
# single run (preserving the non-breaking space before the colon) now
# that nothing (the bookmark) forces them apart any more.
$mergeText = "Les technologies utilisées : html, "
$mergeRange = $d.Content
$mergeFound = $mergeRange.Find.Execute($mergeText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $mergeText, 2)
if (-not $mergeFound) {
    throw "Could not locate 'Les technologies utilisees : html, ' to merge runs"
}

# --- Change 2 ------------------------------------------------------------
$oldTail = "sauf qu’il en reste quelques une que je n’ai pas encore termine comme le fait qu’un joueur ne puisse joueur une même question plus d’une  fois."

# Append a throwaway marker right after the replacement word so the true
# end of the paragraph can be relocated reliably afterwards (adding a
# bookmark exactly at the end of the document's content is unreliable).
$newTail = "demandéesZZMARKERZZ"

$tailRange = $d.Content
$tailFound = $tailRange.Find.Execute($oldTail, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newTail, 2)
if (-not $tailFound) {
    throw "Could not locate the trailing clause to replace with 'demandees'"
}

# Force a run split right before "demandees" (without touching
# formatting) by briefly planting and then removing a bookmark there.
$wordRange = $d.Content
$wordFound = $wordRange.Find.Execute("demandées", $true)
if (-not $wordFound) {
    throw "Could not re-locate 'demandees' after the replacement"
}
$splitPoint = $d.Range($wordRange.Start, $wordRange.Start)
$d.Bookmarks.Add("TEMPSPLIT", $splitPoint)
$d.Bookmarks("TEMPSPLIT").Delete()

# Locate the temporary marker, plant _GoBack right before it (i.e. right
# after "demandees"), then delete the marker text itself.
$markerRange = $d.Content
$markerFound = $markerRange.Find.Execute("ZZMARKERZZ", $true)
if (-not $markerFound) {
    throw "Could not locate the temporary end-of-paragraph marker"
}
$bmPoint = $d.Range($markerRange.Start, $markerRange.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)
$markerRange.Delete()

Write-Host "merge runs found: " $mergeFound
Write-Host "tail replace found: " $tailFound
Write-Host "GoBack bookmark present: " $d.Bookmarks.Exists("_GoBack")
